$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New register (0x0C / row 20) holds the psec4a status register:
$ws.Range("C20").Value = "psec4a_status_register"

# Its sub-fields (rows 18-35, register addresses x0A-x1B) are read-only,
# matching the rest of that block; drop the stray border flag those blank
# cells carried so the formatting is clean like the surrounding read_only
# cells:
$ws.Range("D18:D35").Value = "read_only"
$ws.Range("D18:D35").Borders.LineStyle = -4142

# Document the new ping-pong readout mode bit at the bottom of the
# register map (row 85):
$ws.Range("C85").Value = "psec4a readout mode"
$ws.Range("D85").Value = "0=all samples, 1=ping-pong 528-sample blocks"

# Scroll position / selection moved as part of the edit:
$excel.ActiveWindow.ScrollRow = 76
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D91").Select()
